$wb = $excel.ActiveWorkbook

# --- "Range Status" sheet: species (no.) counts reset to 0, perc. column cleared ---
$wsRange = $wb.Worksheets.Item("Range Status")
$wsRange.Range("B2:B7").Value = 0
$wsRange.Range("C2:C7").ClearContents()

# --- "High Priority break-up" sheet: re-ran classify+summarise, fewer/changed rows ---
$wsBreakup = $wb.Worksheets.Item("High Priority break-up")

# Row 2 stays "Trend New" but with updated counts
$wsBreakup.Range("B2").Value = 12
$wsBreakup.Range("C2").Value = 48
$wsBreakup.Range("D2").Value = 12
$wsBreakup.Range("E2").Value = 48

# Row 3 becomes "IUCN" with its own counts
$wsBreakup.Range("A3").Value = "IUCN"
$wsBreakup.Range("B3").Value = 13
$wsBreakup.Range("C3").Value = 52
$wsBreakup.Range("D3").Value = 13
$wsBreakup.Range("E3").Value = 52

# Old rows 4 ("Range") and 5 ("IUCN") are gone entirely - shift remaining rows up
$wsBreakup.Rows("4:5").Delete()
